$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update sheet (tab) name to reflect the new "through" date
$ws.Name = "Through 2022-09-06"

# 2) Update the month-header label cell (B1) with the new "through" date text
$ws.Cells.Item(1, 2).Value = "September 2022 (through September 06)"

# 3) Rows 5 and 6 swap which neighborhood they represent (Chicago Lawn now
#    sorts above Englewood) and both pick up new 2022-09-06 data, so we
#    rewrite each row fully (column A label plus columns B..BU) to its
#    final target values.
$ws.Cells.Item(5, 1).Value = "Chicago Lawn"
$ws.Cells.Item(6, 1).Value = "Englewood"
$row5Values = @(
    "", 8, "", "", 2, 5, 4, 5, 3, 1, 2, 3, "", "", "", 1, "", 6, "", 3, 2, 2, 1, "", 1, 1, 3, "", 1, "", "", "", 2, 1, 2, 1, "", 1, 1, 1, 3, "", 3, 2, 4, "", 2, 2, 1, 2, "", 1, 2, "", "", 4, "", 2, 2, "", 2, 4, "", "", "", "", "", "", "", 1, "", ""
)

$row6Values = @(
    2, 8, 9, 9, 13, 10, 8, 4, 13, "", 3, 6, 7, 7, 3, 2, "", 4, 1, 2, 9, 8, 7, 2, 4, 3, 5, "", 3, 2, 4, 2, 1, "", 2, 5, 1, 3, 2, 1, 1, 1, 7, 1, 2, 1, 2, 1, 1, "", 4, 1, 3, "", 1, 1, 5, 4, 2, 4, 4, 2, 5, 1, 3, 3, 2, 1, 2, 1, "", 3
)

for ($i = 0; $i -lt $row5Values.Length; $i++) {
    $col = 2 + $i
    $v = $row5Values[$i]
    if ($v -eq "") {
        $ws.Cells.Item(5, $col).ClearContents()
    } else {
        $ws.Cells.Item(5, $col).Value = $v
    }
}

for ($i = 0; $i -lt $row6Values.Length; $i++) {
    $col = 2 + $i
    $v = $row6Values[$i]
    if ($v -eq "") {
        $ws.Cells.Item(6, $col).ClearContents()
    } else {
        $ws.Cells.Item(6, $col).Value = $v
    }
}

# 4) Remaining scattered single-cell updates for 2022-09-06 incidents
#    across other neighborhood rows (new counts added/incremented).
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 20).Value = 2
$ws.Cells.Item(3, 38).Value = 3
$ws.Cells.Item(3, 65).Value = 1
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(7, 20).Value = 1
$ws.Cells.Item(8, 2).Value = 1
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(9, 56).Value = 2
$ws.Cells.Item(10, 29).Value = 1
$ws.Cells.Item(12, 38).Value = 2
$ws.Cells.Item(38, 29).Value = 1
$ws.Cells.Item(43, 2).Value = 1
$ws.Cells.Item(50, 2).Value = 1
$ws.Cells.Item(50, 20).Value = 2
$ws.Cells.Item(77, 65).Value = 1
